# Aggiornamento a l 23 agosto 2021
# Extend the daily-report table on Sheet1 with 14 more days
# (2021-08-10 .. 2021-08-23 / Excel serials 44418..44431),
# carrying the same formatting as the prior last row (343) down
# through the new rows (344..357), with B/C/D = 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 343
$newLastRow = 357

# Copy the formatting of the last existing row down onto the new rows.
$ws.Range("A$lastRow`:D$lastRow").Copy()
$ws.Range("A$($lastRow+1)`:D$newLastRow").PasteSpecial(-4122)

# Fill in the date serials (column A) and zero values (B, C, D).
$startSerial = 44418
for ($r = $lastRow + 1; $r -le $newLastRow; $r++) {
    $serial = $startSerial + ($r - ($lastRow + 1))
    $ws.Cells.Item($r, 1).Value = $serial
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
}
